$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 24, shifting rows 24-27 down to 25-28
$ws.Rows.Item(24).Insert()

# Copy the style of the date cell (D25, formerly D24) into the new D24 cell
$ws.Range("D24").Value = 44505
$ws.Cells.Item(24, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat

$ws.Range("A24").Value = 8
$ws.Range("B24").Value = "Terminal La Palmera de La Serena"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("E24").Value = 4
$ws.Range("F24").Value = 100114007
$ws.Range("G24").Value = "Jengibre"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 16000
$ws.Range("L24").Value = 17000
$ws.Range("M24").Value = 16500
$ws.Range("N24").Value = "$/caja 13 kilos"
$ws.Range("O24").Value = "Perú"
$ws.Range("P24").Value = 1269
$ws.Range("Q24").Value = 13
$ws.Range("R24").Value = "Hortaliza"
